$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 42, pushing nothing (it's currently the
# last row, so this just appends a blank row 42 below the existing row 41).
$ws.Rows.Item(42).Insert()

# Duplicate row 41 (the original last annotation row) down into row 42 via
# copy/paste so cell types (e.g. the text "3" in column B) are preserved
# exactly, without introducing any new cell styles.
$ws.Range("A41:H41").Copy()
$ws.Range("A42").PasteSpecial(-4163)

# Row 42 keeps Annotator/score/expression/purpose (A-D) from the duplicated
# row, but gets its own issue_type, id, source_file and text (E-H).
$ws.Cells.Item(42, 5).Value = "RES"
$ws.Cells.Item(42, 6).Value = "42b1e2ab-785d-481e-b197-1cf6913a8b3e"
$ws.Cells.Item(42, 7).Value = "SJQO7UJCW_annotated.xlsx"
$ws.Cells.Item(42, 8).Value = "However, our main point of the paper is to demonstrate the effectiveness of proposed method against our baseline model shown in Table 1 and 2."

# The original row 41's politeness_score was stored as text "3"; it now
# becomes a proper numeric 3 (matching every other row in the column).
$ws.Cells.Item(41, 2).Value = 3
